$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the price (D) cells that change as Text so the numeric-looking
# strings (e.g. "29.238.52", "0.9990") are preserved exactly, matching
# the original inline-string cell contents instead of being coerced to numbers.
$ws.Range("D2:D16").NumberFormat = "@"
$ws.Range("D18:D29").NumberFormat = "@"
$ws.Range("D31:D31").NumberFormat = "@"
$ws.Range("D33:D50").NumberFormat = "@"

$ws.Range("D2").Value2 = "29.238.52"
$ws.Range("E2").Value2 = "  -0.77%  "
$ws.Range("D3").Value2 = "1.862.55"
$ws.Range("E3").Value2 = "  -0.91%  "
$ws.Range("D4").Value2 = "0.9990"
$ws.Range("E4").Value2 = "  -0.10%  "
$ws.Range("D5").Value2 = "0.7133"
$ws.Range("E5").Value2 = "  -0.69%  "
$ws.Range("D6").Value2 = "240.80"
$ws.Range("E6").Value2 = "  +0.18%  "
$ws.Range("D7").Value2 = "0.9998"
$ws.Range("E7").Value2 = "  -0.02%  "
$ws.Range("D8").Value2 = "0.07736"
$ws.Range("E8").Value2 = "  -1.16%  "
$ws.Range("D9").Value2 = "0.3085"
$ws.Range("E9").Value2 = "  -1.21%  "
$ws.Range("D10").Value2 = "24.97"
$ws.Range("E10").Value2 = "  -0.75%  "
$ws.Range("D11").Value2 = "0.08332"
$ws.Range("E11").Value2 = "  +1.04%  "
$ws.Range("D12").Value2 = "1.876.28"
$ws.Range("E12").Value2 = "  -0.70%  "
$ws.Range("D13").Value2 = "5.204"
$ws.Range("E13").Value2 = "  -1.68%  "
$ws.Range("D14").Value2 = "0.7140"
$ws.Range("E14").Value2 = "  -2.04%  "
$ws.Range("D15").Value2 = "91.20"
$ws.Range("E15").Value2 = "  -0.08%  "
$ws.Range("D16").Value2 = "29.251.43"
$ws.Range("E16").Value2 = "  -1.40%  "
$ws.Range("E17").Value2 = "  -0.05%  "
$ws.Range("D18").Value2 = "242.85"
$ws.Range("E18").Value2 = "  -2.38%  "
$ws.Range("D19").Value2 = "0.000007836"
$ws.Range("E19").Value2 = "  -0.57%  "
$ws.Range("D20").Value2 = "2.127.01"
$ws.Range("E20").Value2 = "  -0.36%  "
$ws.Range("D21").Value2 = "13.18"
$ws.Range("E21").Value2 = "  -1.02%  "
$ws.Range("D22").Value2 = "0.9994"
$ws.Range("E22").Value2 = "  -0.01%  "
$ws.Range("D23").Value2 = "7.898"
$ws.Range("E23").Value2 = "  -1.32%  "
$ws.Range("D24").Value2 = "0.9991"
$ws.Range("E24").Value2 = "  -0.05%  "
$ws.Range("D25").Value2 = "0.1598"
$ws.Range("E25").Value2 = "  +1.96%  "
$ws.Range("D26").Value2 = "163.23"
$ws.Range("D27").Value2 = "8.895"
$ws.Range("E27").Value2 = "  -1.78%  "
$ws.Range("D28").Value2 = "18.50"
$ws.Range("E28").Value2 = "  +0.83%  "
$ws.Range("D29").Value2 = "1.343"
$ws.Range("E29").Value2 = "  -1.61%  "
$ws.Range("E30").Value2 = "  +0.93%  "
$ws.Range("D31").Value2 = "4.421"
$ws.Range("E31").Value2 = "  +0.82%  "
$ws.Range("E32").Value2 = "  +2.42%  "
$ws.Range("D33").Value2 = "0.8362"
$ws.Range("E33").Value2 = "  +15.46%  "
$ws.Range("D34").Value2 = "0.05150"
$ws.Range("E34").Value2 = "  -2.40%  "
$ws.Range("D35").Value2 = "1.932"
$ws.Range("E35").Value2 = "  -0.65%  "
$ws.Range("D36").Value2 = "1.171"
$ws.Range("E36").Value2 = "  -2.72%  "
$ws.Range("D37").Value2 = "2.677"
$ws.Range("E37").Value2 = "  +0.07%  "
$ws.Range("D38").Value2 = "0.01854"
$ws.Range("E38").Value2 = "  -0.53%  "
$ws.Range("D39").Value2 = "2.691"
$ws.Range("E39").Value2 = "  -1.25%  "
$ws.Range("D40").Value2 = "1.174.13"
$ws.Range("E40").Value2 = "  -5.38%  "
$ws.Range("D41").Value2 = "6.192"
$ws.Range("E41").Value2 = "  +1.29%  "
$ws.Range("D42").Value2 = "0.9016"
$ws.Range("E42").Value2 = "  -0.46%  "
$ws.Range("D43").Value2 = "72.82"
$ws.Range("E43").Value2 = "  -1.32%  "
$ws.Range("D44").Value2 = "0.9989"
$ws.Range("E44").Value2 = "  -0.10%  "
$ws.Range("D45").Value2 = "102.30"
$ws.Range("E45").Value2 = "  -1.58%  "
$ws.Range("D46").Value2 = "2.026.10"
$ws.Range("E46").Value2 = "  -0.78%  "
$ws.Range("D47").Value2 = "0.5191"
$ws.Range("E47").Value2 = "  -2.69%  "
$ws.Range("D48").Value2 = "1.790"
$ws.Range("E48").Value2 = "  +1.26%  "
$ws.Range("D49").Value2 = "9.337"
$ws.Range("E49").Value2 = "  +0.42%  "
$ws.Range("D50").Value2 = "7.052"
$ws.Range("E50").Value2 = "  -0.48%  "
$ws.Range("E51").Value2 = "  -0.03%  "
